# Applies the cell-level updates to the cryptos price/volume table.
# D-column numeric-looking values are written with a leading apostrophe
# (Excel "text prefix") so COM keeps them as Text cells (matching the
# original inlineStr cells) instead of auto-converting to Number.
# Style is reset to "Normal" afterward so no stray number-format/style
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.424.09"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.617.93"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'213.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'19.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").Value = "'0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.844.62"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "1.619.62"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "'63.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'236.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.04%  "
$ws.Range("D18").Value = "26.433.93"
$ws.Range("D19").Value = "'7.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.56%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'9.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").Value = "'147.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'7.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D32").Value = "1.526.97"
$ws.Range("E32").Value = "  +6.96%  "
$ws.Range("D33").Value = "'3.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "'1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'0.567"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'5.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.830"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").Value = "1.756.03"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").Value = "'0.762"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'61.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'0.909"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "'90.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.30%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0502"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0961"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.07%  "
